# "Removed Type from imports"
#
# Column J ("Type" header, values "Pool") is dropped from the template.
# The former column K ("Rule For" header, values "Accounting"/"Reporting",
# plus its data-validation dropdown) slides left to become the new J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the text of the "Rule For" header-note living on K1 before anything
# moves, so we can re-home it on the new J1 after the column shift (this
# engine does not re-anchor cell comments automatically on column delete).
$ruleForNote = $ws.Range("K1").Comment.Text()

# Drop the two comments that sit in the column being removed/shifted; the
# "Type" note on J1 goes away for good, the "Rule For" note on K1 is
# recreated at J1 below.
$ws.Range("J1").Comment.Delete()
$ws.Range("K1").Comment.Delete()

# Delete the entire "Type" column (J). Everything to its right - including
# the "Rule For" column and its "Accounting,Reporting" validation list -
# shifts one column to the left automatically.
$ws.Columns.Item(10).Delete()

$ws.Range("J1").AddComment($ruleForNote)

# Match the saved view/selection state (whole-column selection anchored at J1).
$ws.Range("J1:J1048576").Select()
